$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry updates one cell in the cryptos list (Price = D, Volume(1h) = E).
# NumberFormat "@" + resetting Style back to "Normal" afterwards forces the
# assigned value to stay a text string (matching the source inlineStr cells)
# even when the text looks numeric (e.g. "213.23", "88.00"), while leaving the
# cell's effective style/format unchanged.
$updates = @(
    @{ Cell = "D2"; Value = '27.385.47' },
    @{ Cell = "E2"; Value = '  -0.91%  ' },
    @{ Cell = "D3"; Value = '1.654.66' },
    @{ Cell = "E3"; Value = '  -0.16%  ' },
    @{ Cell = "E4"; Value = '  -0.37%  ' },
    @{ Cell = "D5"; Value = '213.23' },
    @{ Cell = "E5"; Value = '  -0.60%  ' },
    @{ Cell = "E6"; Value = '  +0.00%  ' },
    @{ Cell = "E7"; Value = '  -0.38%  ' },
    @{ Cell = "E8"; Value = '  +0.81%  ' },
    @{ Cell = "E9"; Value = '  -0.41%  ' },
    @{ Cell = "E10"; Value = '  -0.95%  ' },
    @{ Cell = "E11"; Value = '  -0.54%  ' },
    @{ Cell = "D12"; Value = '1.889.54' },
    @{ Cell = "E12"; Value = '  -0.15%  ' },
    @{ Cell = "D13"; Value = '1.656.70' },
    @{ Cell = "E13"; Value = '  -0.19%  ' },
    @{ Cell = "E14"; Value = '  -0.89%  ' },
    @{ Cell = "D15"; Value = '0.569' },
    @{ Cell = "E15"; Value = '  +4.00%  ' },
    @{ Cell = "E16"; Value = '  -0.42%  ' },
    @{ Cell = "D17"; Value = '27.383.45' },
    @{ Cell = "E17"; Value = '  -0.87%  ' },
    @{ Cell = "D18"; Value = '231.90' },
    @{ Cell = "E18"; Value = '  -6.02%  ' },
    @{ Cell = "D19"; Value = '0.0₃0727' },
    @{ Cell = "E19"; Value = '  -0.17%  ' },
    @{ Cell = "E20"; Value = '  -0.17%  ' },
    @{ Cell = "E21"; Value = '  -0.26%  ' },
    @{ Cell = "E22"; Value = '  -2.10%  ' },
    @{ Cell = "D23"; Value = '9.40' },
    @{ Cell = "E23"; Value = '  +3.60%  ' },
    @{ Cell = "E24"; Value = '  +0.10%  ' },
    @{ Cell = "D25"; Value = '147.42' },
    @{ Cell = "E25"; Value = '  +0.68%  ' },
    @{ Cell = "D26"; Value = '7.10' },
    @{ Cell = "E26"; Value = '  -0.80%  ' },
    @{ Cell = "E27"; Value = '  -1.92%  ' },
    @{ Cell = "E28"; Value = '  -0.39%  ' },
    @{ Cell = "E29"; Value = '  +0.27%  ' },
    @{ Cell = "E30"; Value = '  -0.36%  ' },
    @{ Cell = "E31"; Value = '  -4.24%  ' },
    @{ Cell = "E32"; Value = '  -1.30%  ' },
    @{ Cell = "E33"; Value = '  +0.14%  ' },
    @{ Cell = "D34"; Value = '1.422.15' },
    @{ Cell = "E34"; Value = '  -0.46%  ' },
    @{ Cell = "E35"; Value = '  +1.30%  ' },
    @{ Cell = "E36"; Value = '  -1.02%  ' },
    @{ Cell = "D37"; Value = '0.907' },
    @{ Cell = "E37"; Value = '  -2.01%  ' },
    @{ Cell = "D38"; Value = '0.571' },
    @{ Cell = "E38"; Value = '  -1.26%  ' },
    @{ Cell = "E39"; Value = '  +0.09%  ' },
    @{ Cell = "D40"; Value = '1.05' },
    @{ Cell = "E40"; Value = '  +0.10%  ' },
    @{ Cell = "E41"; Value = '  -0.44%  ' },
    @{ Cell = "D42"; Value = '5.52' },
    @{ Cell = "E42"; Value = '  +2.45%  ' },
    @{ Cell = "D43"; Value = '0.798' },
    @{ Cell = "E43"; Value = '  +1.23%  ' },
    @{ Cell = "E44"; Value = '  +0.25%  ' },
    @{ Cell = "D45"; Value = '64.93' },
    @{ Cell = "E45"; Value = '  -6.04%  ' },
    @{ Cell = "D46"; Value = '1.797.93' },
    @{ Cell = "E46"; Value = '  -0.06%  ' },
    @{ Cell = "E47"; Value = '  -1.29%  ' },
    @{ Cell = "D48"; Value = '88.00' },
    @{ Cell = "E48"; Value = '  -0.67%  ' },
    @{ Cell = "E49"; Value = '  +1.57%  ' },
    @{ Cell = "E50"; Value = '  -0.01%  ' },
    @{ Cell = "E51"; Value = '  -0.48%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).NumberFormat = "@"
    $ws.Range($u.Cell).Value = $u.Value
    $ws.Range($u.Cell).Style = "Normal"
}
